# "Results section done. Next up: Discussion"
#
# Adds the new column-I figures alongside the DP54 block (rows 47-50) and
# replaces the stray leftover row 53 with a proper new data row (row 54)
# that follows the same tol/calls/eval pattern as the rows above it.
# Finally restores the scrolled view with I51 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I figures next to the DP54 block (rows 47-50) -----------
$ws.Range("I47").Formula = "=911.915-15.83/2"
$ws.Range("I48").Formula = "=407.915+292.33"
$ws.Range("I49").Formula = "=I48-15.383/2"
$ws.Range("I50").Formula = "=38.755+15.383/2"

# --- Replace the stray row 53 with a proper new data row 54 -------------
# Row 53 previously held just one leftover cell (E53 = 1.333*E52). Remove
# the entire row and add a fresh data row one line further down, matching
# the tol/calls/eval pattern used by rows 43-52 (tol=1e-12, calls=345,
# eval=13*D54).
$null = $ws.Rows("53").Delete()

$ws.Range("C54").Value = 0.000000000001
$ws.Range("D54").Value = 345
$ws.Range("E54").Formula = "=13*D54"

# --- Restore the view: scrolled down with I51 selected -------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("I51").Select()
